$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 52: new timesheet entry
$ws.Cells.Item(52, 1).Value = 44055
$ws.Cells.Item(52, 2).Value = 2
$ws.Cells.Item(52, 3).Value = "przeniesienie modelu do osobnego projektu. Dłubanie nad serwerem iis, niestety "

# Row 53: new timesheet entry (wrapped comment cell, matching the sheet's alternating-row style)
$ws.Cells.Item(53, 1).Value = 44055
$ws.Cells.Item(53, 2).Value = 2
$ws.Cells.Item(53, 3).Value = "firebase"
$ws.Cells.Item(53, 3).WrapText = $true

# Move the selection down to the next empty comment cell, like the author did
$ws.Range("C54").Select()
